# This script applies a reshuffle of the weekly Fruta/Hortaliza price rows
# (rows 4-10) in the active worksheet. The underlying commit re-ordered the
# daily price observations (each row is one market day), so several rows
# end up swapping their Fecha / Volumen / Precio / Origen values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each row, taken from the row that supplied the data after
# the reshuffle (row -> source row): 4<-9, 5<-7, 6<-8, 7<-5, 8<-6, 9<-10, 10<-4
$updates = @(
    @{ Row = 4;  D = 44540; M = 240; N = 3500; O = 3800; P = 3650; R = "Región del Maule";       S = 1825 },
    @{ Row = 5;  D = 44181; M = 65;  N = 3600; O = 3800; P = 3692; R = "Provincia de Diguillín";  S = 1846 },
    @{ Row = 6;  D = 44181; M = 80;  N = 1800; O = 2000; P = 1875;                                S = 1875 },
    @{ Row = 7;  D = 44187; M = 80;  N = 2800; O = 3000; P = 2900; R = "Provincia de Linares";    S = 1450 },
    @{ Row = 8;  D = 44187; M = 65;  N = 1400; O = 1500; P = 1446;                                S = 1446 },
    @{ Row = 9;  D = 44596; M = 120; N = 2500; O = 2700; P = 2600; R = "Provincia de Linares";    S = 1300 },
    @{ Row = 10; D = 44594;                    O = 2800; P = 2650;                                S = 1325 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value = $u.D    # D: Fecha
    if ($u.ContainsKey("M")) { $ws.Cells.Item($r, 13).Value = $u.M }   # M: Volumen
    if ($u.ContainsKey("N")) { $ws.Cells.Item($r, 14).Value = $u.N }   # N: Precio mínimo
    $ws.Cells.Item($r, 15).Value = $u.O   # O: Precio máximo
    $ws.Cells.Item($r, 16).Value = $u.P   # P: Precio promedio ponderado
    if ($u.ContainsKey("R")) { $ws.Cells.Item($r, 18).Value = $u.R }   # R: Origen
    $ws.Cells.Item($r, 19).Value = $u.S   # S: Precio $/Kg
}
